$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 612 (pushes the existing rows 612-638 down to 613-639)
$ws.Rows.Item(612).Insert()

# Populate the newly inserted row with this week's price record for
# Feria Lagunitas de Puerto Montt - Coliflor (same metadata as the
# surrounding rows, new date + price figures)
$ws.Cells.Item(612, 1).Value = 4
$ws.Cells.Item(612, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(612, 3).Value = "Los Lagos"
$ws.Cells.Item(612, 4).Value = 45267
$ws.Cells.Item(612, 5).Value = 10
$ws.Cells.Item(612, 6).Value = 100112008
$ws.Cells.Item(612, 7).Value = "Coliflor"
$ws.Cells.Item(612, 8).Value = "Sin especificar"
$ws.Cells.Item(612, 9).Value = "Primera"
$ws.Cells.Item(612, 10).Value = 500
$ws.Cells.Item(612, 11).Value = 1600
$ws.Cells.Item(612, 12).Value = 1600
$ws.Cells.Item(612, 13).Value = 1600
$ws.Cells.Item(612, 14).Value = "`$/unidad"
$ws.Cells.Item(612, 15).Value = "Región Metropolitana"
$ws.Cells.Item(612, 16).Value = 1600
$ws.Cells.Item(612, 17).Value = 1
$ws.Cells.Item(612, 18).Value = "Hortaliza"
